# Fruta / hortaliza, semanal
# Insert a new weekly pricing observation (two rows: Primera / Segunda
# quality) at the top of the Betarraga data block. The existing rows
# shift down by two positions to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 406:425 down to 408:427, inserting two fresh blank rows at 406:407.
$ws.Rows("406:407").Insert()

# Row 406: new "Primera" quality observation.
$ws.Cells.Item(406, 1).Value = 8
$ws.Cells.Item(406, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(406, 3).Value = "Coquimbo"
$ws.Cells.Item(406, 4).Value = 45008
$ws.Cells.Item(406, 5).Value = 4
$ws.Cells.Item(406, 6).Value = 100114014
$ws.Cells.Item(406, 7).Value = "Betarraga"
$ws.Cells.Item(406, 8).Value = "Sin especificar"
$ws.Cells.Item(406, 9).Value = "Primera"
$ws.Cells.Item(406, 10).Value = 2000
$ws.Cells.Item(406, 11).Value = 500
$ws.Cells.Item(406, 12).Value = 600
$ws.Cells.Item(406, 13).Value = 550
$ws.Cells.Item(406, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(406, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(406, 16).Value = 183
$ws.Cells.Item(406, 17).Value = 3
$ws.Cells.Item(406, 18).Value = "Hortaliza"

# Row 407: new "Segunda" quality observation, same date.
$ws.Cells.Item(407, 1).Value = 8
$ws.Cells.Item(407, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(407, 3).Value = "Coquimbo"
$ws.Cells.Item(407, 4).Value = 45008
$ws.Cells.Item(407, 5).Value = 4
$ws.Cells.Item(407, 6).Value = 100114014
$ws.Cells.Item(407, 7).Value = "Betarraga"
$ws.Cells.Item(407, 8).Value = "Sin especificar"
$ws.Cells.Item(407, 9).Value = "Segunda"
$ws.Cells.Item(407, 10).Value = 1500
$ws.Cells.Item(407, 11).Value = 400
$ws.Cells.Item(407, 12).Value = 450
$ws.Cells.Item(407, 13).Value = 425
$ws.Cells.Item(407, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(407, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(407, 16).Value = 142
$ws.Cells.Item(407, 17).Value = 3
$ws.Cells.Item(407, 18).Value = "Hortaliza"
